$d = $word.ActiveDocument

$d.Content.Find.Execute("47÷8=5, 7", $true, $false, $false, $false, $false, $true, 1, $false, "15÷2=7, 1", 2) | Out-Null
$d.Content.Find.Execute("43÷2=21, 1", $true, $false, $false, $false, $false, $true, 1, $false, "91÷9=10, 1", 2) | Out-Null
$d.Content.Find.Execute("84÷5=16, 4", $true, $false, $false, $false, $false, $true, 1, $false, "31÷4=7, 3", 2) | Out-Null
$d.Content.Find.Execute("18÷8=2, 2", $true, $false, $false, $false, $false, $true, 1, $false, "33÷5=6, 3", 2) | Out-Null
$d.Content.Find.Execute("71÷3=23, 2", $true, $false, $false, $false, $false, $true, 1, $false, "59÷9=6, 5", 2) | Out-Null
$d.Content.Find.Execute("65÷2=32, 1", $true, $false, $false, $false, $false, $true, 1, $false, "54÷9=6, 0", 2) | Out-Null
$d.Content.Find.Execute("41÷4=10, 1", $true, $false, $false, $false, $false, $true, 1, $false, "72÷9=8, 0", 2) | Out-Null
$d.Content.Find.Execute("64÷4=16, 0", $true, $false, $false, $false, $false, $true, 1, $false, "73÷4=18, 1", 2) | Out-Null
$d.Content.Find.Execute("54÷8=6, 6", $true, $false, $false, $false, $false, $true, 1, $false, "99÷7=14, 1", 2) | Out-Null
$d.Content.Find.Execute("32÷6=5, 2", $true, $false, $false, $false, $false, $true, 1, $false, "24÷9=2, 6", 2) | Out-Null
$d.Content.Find.Execute("70÷4=17, 2", $true, $false, $false, $false, $false, $true, 1, $false, "63÷8=7, 7", 2) | Out-Null
$d.Content.Find.Execute("96÷9=10, 6", $true, $false, $false, $false, $false, $true, 1, $false, "56÷9=6, 2", 2) | Out-Null
$d.Content.Find.Execute("73÷7=10, 3", $true, $false, $false, $false, $false, $true, 1, $false, "47÷9=5, 2", 2) | Out-Null
$d.Content.Find.Execute("61÷2=30, 1", $true, $false, $false, $false, $false, $true, 1, $false, "76÷8=9, 4", 2) | Out-Null
$d.Content.Find.Execute("10÷2=5, 0", $true, $false, $false, $false, $false, $true, 1, $false, "31÷8=3, 7", 2) | Out-Null
$d.Content.Find.Execute("42÷5=8, 2", $true, $false, $false, $false, $false, $true, 1, $false, "15÷5=3, 0", 2) | Out-Null
$d.Content.Find.Execute("34÷3=11, 1", $true, $false, $false, $false, $false, $true, 1, $false, "52÷6=8, 4", 2) | Out-Null
$d.Content.Find.Execute("21÷9=2, 3", $true, $false, $false, $false, $false, $true, 1, $false, "99÷8=12, 3", 2) | Out-Null
$d.Content.Find.Execute("37÷7=5, 2", $true, $false, $false, $false, $false, $true, 1, $false, "84÷8=10, 4", 2) | Out-Null
$d.Content.Find.Execute("90÷5=18, 0", $true, $false, $false, $false, $false, $true, 1, $false, "43÷3=14, 1", 2) | Out-Null
$d.Content.Find.Execute("60÷3=20, 0", $true, $false, $false, $false, $false, $true, 1, $false, "12÷9=1, 3", 2) | Out-Null
$d.Content.Find.Execute("50÷2=25, 0", $true, $false, $false, $false, $false, $true, 1, $false, "17÷2=8, 1", 2) | Out-Null
$d.Content.Find.Execute("62÷6=10, 2", $true, $false, $false, $false, $false, $true, 1, $false, "54÷9=6, 0", 2) | Out-Null
$d.Content.Find.Execute("60÷7=8, 4", $true, $false, $false, $false, $false, $true, 1, $false, "60÷6=10, 0", 2) | Out-Null
$d.Content.Find.Execute("52÷2=26, 0", $true, $false, $false, $false, $false, $true, 1, $false, "19÷9=2, 1", 2) | Out-Null
